# "Generate Report for handoff": a new handoff just occurred for the
# 5a562851-...-.md source file, so its "Latest Handoff Datetime" (column D)
# is refreshed with a newer timestamp, for both locale report sheets.
$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("D4").Value = "2016-01-26 09:17:12"
$wsDe.Range("D4").Value = "2016-01-26 09:17:23"
